$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 0.5021227866937743
$ws.Range("C5").Value = 0.03660012403874354
$ws.Range("D5").Value = 0.09310312614494011

$ws.Range("B7").Value = 5.950933613741086
$ws.Range("C7").Value = 0.17589663884708953
$ws.Range("D7").Value = 1.8810431795225795

$ws.Range("B8").Value = 0.5397717538274582
$ws.Range("C8").Value = 0.045135982952844225
$ws.Range("D8").Value = 0.11176427719582135

$ws.Range("B11").Value = 1.7618568524032712
$ws.Range("C11").Value = 0.03115566246519964
$ws.Range("D11").Value = 0.2965747016567983

$ws.Range("B12").Value = 4.02971547044312
$ws.Range("C12").Value = 0.12578524359393609
$ws.Range("D12").Value = 2.212279906641036

$ws.Range("B14").Value = 2.992671955117957
$ws.Range("C14").Value = 0.08481251966043689
$ws.Range("D14").Value = 0.8132714584011018

$ws.Range("B17").Value = 0.423819427194849
$ws.Range("C17").Value = 0.0035923167312428617
$ws.Range("D17").Value = 0.256143598007691

$ws.Range("B18").Value = 1.0676707032386923
$ws.Range("C18").Value = 0.0416588022222932
$ws.Range("D18").Value = 1.0318023415418591

$ws.Range("B19").Value = 1.22576240986613
$ws.Range("C19").Value = 0.08165675410834797
$ws.Range("D19").Value = 0.9070743830604145

$ws.Range("B20").Value = 9.7176908415361
$ws.Range("C20").Value = 0.5044040593754678
$ws.Range("D20").Value = 2.4970016104618855

$ws.Range("B21").Value = 1.9174302782375319
$ws.Range("C21").Value = 0.10215175817893439
$ws.Range("D21").Value = 0.48000277419155385

$ws.Range("B22").Value = 1.131529124430239
$ws.Range("C22").Value = 0.034380954564611516
$ws.Range("D22").Value = 0.41741473811975177

$ws.Range("B23").Value = 1.889234381536828
$ws.Range("C23").Value = 0.07551142596317245
$ws.Range("D23").Value = 1.7548665520775215

$ws.Range("B24").Value = 0.419101513419447
$ws.Range("C24").Value = 0.03245134119506996
$ws.Range("D24").Value = 0.05358389280606204

$ws.Range("B25").Value = 2.315193647365876
$ws.Range("C25").Value = 0.020106100758396684
$ws.Range("D25").Value = 1.6733228076345497

$ws.Range("B26").Value = 1.8711080690234048
$ws.Range("C26").Value = 0.10573056497720043
$ws.Range("D26").Value = 0.5720516929044994

$ws.Range("B27").Value = 3.1404265720573035
$ws.Range("C27").Value = 0.16409575628050913
$ws.Range("D27").Value = 0.8656063689436064

$ws.Range("B28").Value = 2.056460600201815
$ws.Range("C28").Value = 0.09927374238659825
$ws.Range("D28").Value = 0.5008409065290582

